# Remove the trailing "Ver no Jupiter Salvar em pdf Salvar em docx" line,
# the "© 2020 . Contact: ..." footer line, and the blank paragraph that
# separates them from the final page-break paragraph at the very end of
# the document. The blank paragraph right after the bibliography entry
# ("...Editora Ática.") is kept as-is.

$d = $word.ActiveDocument

$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $p
    }
}

$copyrightPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $copyrightPara = $p
    }
}

# The blank paragraph that immediately follows the copyright line (the
# second of the two blank paragraphs preceding the page-break paragraph)
# must also be removed.
$afterCopyright = $copyrightPara.Next()

$deleteRange = $d.Range($startPara.Range.Start, $afterCopyright.Range.End)
$deleteRange.Delete()
